$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 34; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value -eq 45233) {
        $cell.Value = 45243
    }
}
